$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A23").Value = 6
$ws.Range("B23").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C23").Value = "Metropolitana"
$ws.Range("D23").Value = 44533
$ws.Range("E23").Value = 13
$ws.Range("F23").Value = "Fruta"
$ws.Range("G23").Value = 100103
$ws.Range("H23").Value = "Frutos de hueso (carozo)"
$ws.Range("I23").Value = 100103003
$ws.Range("J23").Value = "Damasco"
$ws.Range("K23").Value = "Castle Brite"
$ws.Range("L23").Value = "Primera"
$ws.Range("M23").Value = 12
$ws.Range("N23").Value = 400000
$ws.Range("O23").Value = 400000
$ws.Range("P23").Value = 400000
$ws.Range("Q23").Value = "`$/bins (500 kilos)"
$ws.Range("R23").Value = "Región Metropolitana"
$ws.Range("S23").Value = 800
$ws.Range("T23").Value = 500

$ws.Range("A24").Value = 6
$ws.Range("B24").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C24").Value = "Metropolitana"
$ws.Range("D24").Value = 44533
$ws.Range("E24").Value = 13
$ws.Range("F24").Value = "Fruta"
$ws.Range("G24").Value = 100103
$ws.Range("H24").Value = "Frutos de hueso (carozo)"
$ws.Range("I24").Value = 100103003
$ws.Range("J24").Value = "Damasco"
$ws.Range("K24").Value = "Castle Brite"
$ws.Range("L24").Value = "Segunda"
$ws.Range("M24").Value = 12
$ws.Range("N24").Value = 370000
$ws.Range("O24").Value = 370000
$ws.Range("P24").Value = 370000
$ws.Range("Q24").Value = "`$/bins (500 kilos)"
$ws.Range("R24").Value = "Región Metropolitana"
$ws.Range("S24").Value = 740
$ws.Range("T24").Value = 500

$ws.Range("A25").Value = 6
$ws.Range("B25").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C25").Value = "Metropolitana"
$ws.Range("D25").Value = 44172
$ws.Range("E25").Value = 13
$ws.Range("F25").Value = "Fruta"
$ws.Range("G25").Value = 100103
$ws.Range("H25").Value = "Frutos de hueso (carozo)"
$ws.Range("I25").Value = 100103003
$ws.Range("J25").Value = "Damasco"
$ws.Range("K25").Value = "Castle Brite"
$ws.Range("L25").Value = "Primera"
$ws.Range("M25").Value = 300
$ws.Range("N25").Value = 10000
$ws.Range("O25").Value = 11000
$ws.Range("P25").Value = 10500
$ws.Range("Q25").Value = "`$/caja 10 kilos"
$ws.Range("R25").Value = "Región de O'Higgins"
$ws.Range("S25").Value = 1050
$ws.Range("T25").Value = 10

$ws.Range("A26").Value = 6
$ws.Range("B26").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C26").Value = "Metropolitana"
$ws.Range("D26").Value = 44172
$ws.Range("E26").Value = 13
$ws.Range("F26").Value = "Fruta"
$ws.Range("G26").Value = 100103
$ws.Range("H26").Value = "Frutos de hueso (carozo)"
$ws.Range("I26").Value = 100103003
$ws.Range("J26").Value = "Damasco"
$ws.Range("K26").Value = "Castle Brite"
$ws.Range("L26").Value = "Segunda"
$ws.Range("M26").Value = 100
$ws.Range("N26").Value = 10000
$ws.Range("O26").Value = 10000
$ws.Range("P26").Value = 10000
$ws.Range("Q26").Value = "`$/caja 18 kilos"
$ws.Range("R26").Value = "Región de O'Higgins"
$ws.Range("S26").Value = 556
$ws.Range("T26").Value = 18

$ws.Range("A27").Value = 6
$ws.Range("B27").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C27").Value = "Metropolitana"
$ws.Range("D27").Value = 44186
$ws.Range("E27").Value = 13
$ws.Range("F27").Value = "Fruta"
$ws.Range("G27").Value = 100103
$ws.Range("H27").Value = "Frutos de hueso (carozo)"
$ws.Range("I27").Value = 100103003
$ws.Range("J27").Value = "Damasco"
$ws.Range("K27").Value = "Dina"
$ws.Range("L27").Value = "Especial"
$ws.Range("M27").Value = 250
$ws.Range("N27").Value = 17000
$ws.Range("O27").Value = 17000
$ws.Range("P27").Value = 17000
$ws.Range("Q27").Value = "`$/caja 18 kilos"
$ws.Range("R27").Value = "Región Metropolitana"
$ws.Range("S27").Value = 944
$ws.Range("T27").Value = 18

$ws.Range("A28").Value = 6
$ws.Range("B28").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C28").Value = "Metropolitana"
$ws.Range("D28").Value = 44186
$ws.Range("E28").Value = 13
$ws.Range("F28").Value = "Fruta"
$ws.Range("G28").Value = 100103
$ws.Range("H28").Value = "Frutos de hueso (carozo)"
$ws.Range("I28").Value = 100103003
$ws.Range("J28").Value = "Damasco"
$ws.Range("K28").Value = "Dina"
$ws.Range("L28").Value = "Primera"
$ws.Range("M28").Value = 150
$ws.Range("N28").Value = 17000
$ws.Range("O28").Value = 17000
$ws.Range("P28").Value = 17000
$ws.Range("Q28").Value = "`$/caja 18 kilos empedrada"
$ws.Range("R28").Value = "Región de O'Higgins"
$ws.Range("S28").Value = 944
$ws.Range("T28").Value = 18

$ws.Range("A29").Value = 6
$ws.Range("B29").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C29").Value = "Metropolitana"
$ws.Range("D29").Value = 44179
$ws.Range("E29").Value = 13
$ws.Range("F29").Value = "Fruta"
$ws.Range("G29").Value = 100103
$ws.Range("H29").Value = "Frutos de hueso (carozo)"
$ws.Range("I29").Value = 100103003
$ws.Range("J29").Value = "Damasco"
$ws.Range("K29").Value = "Dina"
$ws.Range("L29").Value = "Primera"
$ws.Range("M29").Value = 120
$ws.Range("N29").Value = 11000
$ws.Range("O29").Value = 11000
$ws.Range("P29").Value = 11000
$ws.Range("Q29").Value = "`$/caja 16 kilos"
$ws.Range("R29").Value = "Región Metropolitana"
$ws.Range("S29").Value = 688
$ws.Range("T29").Value = 16

$ws.Range("A30").Value = 6
$ws.Range("B30").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C30").Value = "Metropolitana"
$ws.Range("D30").Value = 44179
$ws.Range("E30").Value = 13
$ws.Range("F30").Value = "Fruta"
$ws.Range("G30").Value = 100103
$ws.Range("H30").Value = "Frutos de hueso (carozo)"
$ws.Range("I30").Value = 100103003
$ws.Range("J30").Value = "Damasco"
$ws.Range("K30").Value = "Dina"
$ws.Range("L30").Value = "Segunda"
$ws.Range("M30").Value = 170
$ws.Range("N30").Value = 9000
$ws.Range("O30").Value = 9000
$ws.Range("P30").Value = 9000
$ws.Range("Q30").Value = "`$/caja 16 kilos"
$ws.Range("R30").Value = "Región Metropolitana"
$ws.Range("S30").Value = 562
$ws.Range("T30").Value = 16

$ws.Range("A31").Value = 6
$ws.Range("B31").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C31").Value = "Metropolitana"
$ws.Range("D31").Value = 44179
$ws.Range("E31").Value = 13
$ws.Range("F31").Value = "Fruta"
$ws.Range("G31").Value = 100103
$ws.Range("H31").Value = "Frutos de hueso (carozo)"
$ws.Range("I31").Value = 100103003
$ws.Range("J31").Value = "Damasco"
$ws.Range("K31").Value = "Modesto"
$ws.Range("L31").Value = "Primera"
$ws.Range("M31").Value = 253
$ws.Range("N31").Value = 10000
$ws.Range("O31").Value = 12000
$ws.Range("P31").Value = 10988
$ws.Range("Q31").Value = "`$/caja 12 kilos"
$ws.Range("R31").Value = "Región de O'Higgins"
$ws.Range("S31").Value = 916
$ws.Range("T31").Value = 12

$ws.Range("A32").Value = 6
$ws.Range("B32").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C32").Value = "Metropolitana"
$ws.Range("D32").Value = 44179
$ws.Range("E32").Value = 13
$ws.Range("F32").Value = "Fruta"
$ws.Range("G32").Value = 100103
$ws.Range("H32").Value = "Frutos de hueso (carozo)"
$ws.Range("I32").Value = 100103003
$ws.Range("J32").Value = "Damasco"
$ws.Range("K32").Value = "Modesto"
$ws.Range("L32").Value = "Segunda"
$ws.Range("M32").Value = 170
$ws.Range("N32").Value = 8000
$ws.Range("O32").Value = 8000
$ws.Range("P32").Value = 8000
$ws.Range("Q32").Value = "`$/caja 12 kilos"
$ws.Range("R32").Value = "Región de O'Higgins"
$ws.Range("S32").Value = 667
$ws.Range("T32").Value = 12

$ws.Range("A33").Value = 6
$ws.Range("B33").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C33").Value = "Metropolitana"
$ws.Range("D33").Value = 44195
$ws.Range("E33").Value = 13
$ws.Range("F33").Value = "Fruta"
$ws.Range("G33").Value = 100103
$ws.Range("H33").Value = "Frutos de hueso (carozo)"
$ws.Range("I33").Value = 100103003
$ws.Range("J33").Value = "Damasco"
$ws.Range("K33").Value = "Albaricoque"
$ws.Range("L33").Value = "Segunda"
$ws.Range("M33").Value = 200
$ws.Range("N33").Value = 12000
$ws.Range("O33").Value = 12000
$ws.Range("P33").Value = 12000
$ws.Range("Q33").Value = "`$/caja 18 kilos"
$ws.Range("R33").Value = "Región Metropolitana"
$ws.Range("S33").Value = 667
$ws.Range("T33").Value = 18

$ws.Range("A34").Value = 6
$ws.Range("B34").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C34").Value = "Metropolitana"
$ws.Range("D34").Value = 44176
$ws.Range("E34").Value = 13
$ws.Range("F34").Value = "Fruta"
$ws.Range("G34").Value = 100103
$ws.Range("H34").Value = "Frutos de hueso (carozo)"
$ws.Range("I34").Value = 100103003
$ws.Range("J34").Value = "Damasco"
$ws.Range("K34").Value = "Castle Brite"
$ws.Range("L34").Value = "Primera"
$ws.Range("M34").Value = 250
$ws.Range("N34").Value = 9000
$ws.Range("O34").Value = 9000
$ws.Range("P34").Value = 9000
$ws.Range("Q34").Value = "`$/caja 12 kilos"
$ws.Range("R34").Value = "Región Metropolitana"
$ws.Range("S34").Value = 750
$ws.Range("T34").Value = 12

$ws.Range("A35").Value = 6
$ws.Range("B35").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C35").Value = "Metropolitana"
$ws.Range("D35").Value = 44176
$ws.Range("E35").Value = 13
$ws.Range("F35").Value = "Fruta"
$ws.Range("G35").Value = 100103
$ws.Range("H35").Value = "Frutos de hueso (carozo)"
$ws.Range("I35").Value = 100103003
$ws.Range("J35").Value = "Damasco"
$ws.Range("K35").Value = "Castle Brite"
$ws.Range("L35").Value = "Segunda"
$ws.Range("M35").Value = 250
$ws.Range("N35").Value = 8000
$ws.Range("O35").Value = 8000
$ws.Range("P35").Value = 8000
$ws.Range("Q35").Value = "`$/caja 12 kilos"
$ws.Range("R35").Value = "Región Metropolitana"
$ws.Range("S35").Value = 667
$ws.Range("T35").Value = 12

$ws.Range("A36").Value = 6
$ws.Range("B36").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C36").Value = "Metropolitana"
$ws.Range("D36").Value = 44165
$ws.Range("E36").Value = 13
$ws.Range("F36").Value = "Fruta"
$ws.Range("G36").Value = 100103
$ws.Range("H36").Value = "Frutos de hueso (carozo)"
$ws.Range("I36").Value = 100103003
$ws.Range("J36").Value = "Damasco"
$ws.Range("K36").Value = "Castle Brite"
$ws.Range("L36").Value = "Primera"
$ws.Range("M36").Value = 160
$ws.Range("N36").Value = 13000
$ws.Range("O36").Value = 14000
$ws.Range("P36").Value = 13500
$ws.Range("Q36").Value = "`$/bandeja 6 kilos"
$ws.Range("R36").Value = "Región Metropolitana"
$ws.Range("S36").Value = 2250
$ws.Range("T36").Value = 6

$ws.Range("A37").Value = 6
$ws.Range("B37").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C37").Value = "Metropolitana"
$ws.Range("D37").Value = 44530
$ws.Range("E37").Value = 13
$ws.Range("F37").Value = "Fruta"
$ws.Range("G37").Value = 100103
$ws.Range("H37").Value = "Frutos de hueso (carozo)"
$ws.Range("I37").Value = 100103003
$ws.Range("J37").Value = "Damasco"
$ws.Range("K37").Value = "Castle Brite"
$ws.Range("L37").Value = "Primera"
$ws.Range("M37").Value = 14
$ws.Range("N37").Value = 500000
$ws.Range("O37").Value = 550000
$ws.Range("P37").Value = 525000
$ws.Range("Q37").Value = "`$/bins (500 kilos)"
$ws.Range("R37").Value = "Región Metropolitana"
$ws.Range("S37").Value = 1050
$ws.Range("T37").Value = 500

$ws.Range("A38").Value = 6
$ws.Range("B38").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C38").Value = "Metropolitana"
$ws.Range("D38").Value = 44530
$ws.Range("E38").Value = 13
$ws.Range("F38").Value = "Fruta"
$ws.Range("G38").Value = 100103
$ws.Range("H38").Value = "Frutos de hueso (carozo)"
$ws.Range("I38").Value = 100103003
$ws.Range("J38").Value = "Damasco"
$ws.Range("K38").Value = "Castle Brite"
$ws.Range("L38").Value = "Segunda"
$ws.Range("M38").Value = 12
$ws.Range("N38").Value = 400000
$ws.Range("O38").Value = 400000
$ws.Range("P38").Value = 400000
$ws.Range("Q38").Value = "`$/bins (500 kilos)"
$ws.Range("R38").Value = "Región Metropolitana"
$ws.Range("S38").Value = 800
$ws.Range("T38").Value = 500

$ws.Range("A39").Value = 6
$ws.Range("B39").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C39").Value = "Metropolitana"
$ws.Range("D39").Value = 44188
$ws.Range("E39").Value = 13
$ws.Range("F39").Value = "Fruta"
$ws.Range("G39").Value = 100103
$ws.Range("H39").Value = "Frutos de hueso (carozo)"
$ws.Range("I39").Value = 100103003
$ws.Range("J39").Value = "Damasco"
$ws.Range("K39").Value = "Dina"
$ws.Range("L39").Value = "Especial"
$ws.Range("M39").Value = 75
$ws.Range("N39").Value = 18000
$ws.Range("O39").Value = 18000
$ws.Range("P39").Value = 18000
$ws.Range("Q39").Value = "`$/caja 18 kilos"
$ws.Range("R39").Value = "Región de O'Higgins"
$ws.Range("S39").Value = 1000
$ws.Range("T39").Value = 18

$ws.Range("A40").Value = 6
$ws.Range("B40").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C40").Value = "Metropolitana"
$ws.Range("D40").Value = 44188
$ws.Range("E40").Value = 13
$ws.Range("F40").Value = "Fruta"
$ws.Range("G40").Value = 100103
$ws.Range("H40").Value = "Frutos de hueso (carozo)"
$ws.Range("I40").Value = 100103003
$ws.Range("J40").Value = "Damasco"
$ws.Range("K40").Value = "Dina"
$ws.Range("L40").Value = "Primera"
$ws.Range("M40").Value = 170
$ws.Range("N40").Value = 15000
$ws.Range("O40").Value = 16000
$ws.Range("P40").Value = 15500
$ws.Range("Q40").Value = "`$/caja 18 kilos"
$ws.Range("R40").Value = "Región de O'Higgins"
$ws.Range("S40").Value = 861
$ws.Range("T40").Value = 18

$ws.Range("A41").Value = 6
$ws.Range("B41").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C41").Value = "Metropolitana"
$ws.Range("D41").Value = 44159
$ws.Range("E41").Value = 13
$ws.Range("F41").Value = "Fruta"
$ws.Range("G41").Value = 100103
$ws.Range("H41").Value = "Frutos de hueso (carozo)"
$ws.Range("I41").Value = 100103003
$ws.Range("J41").Value = "Damasco"
$ws.Range("K41").Value = "Castle Brite"
$ws.Range("L41").Value = "Primera"
$ws.Range("M41").Value = 170
$ws.Range("N41").Value = 14000
$ws.Range("O41").Value = 15000
$ws.Range("P41").Value = 14500
$ws.Range("Q41").Value = "`$/caja 15 kilos"
$ws.Range("R41").Value = "Región Metropolitana"
$ws.Range("S41").Value = 967
$ws.Range("T41").Value = 15

$ws.Range("A42").Value = 6
$ws.Range("B42").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C42").Value = "Metropolitana"
$ws.Range("D42").Value = 44519
$ws.Range("E42").Value = 13
$ws.Range("F42").Value = "Fruta"
$ws.Range("G42").Value = 100103
$ws.Range("H42").Value = "Frutos de hueso (carozo)"
$ws.Range("I42").Value = 100103003
$ws.Range("J42").Value = "Damasco"
$ws.Range("K42").Value = "Castle Brite"
$ws.Range("L42").Value = "Primera"
$ws.Range("M42").Value = 150
$ws.Range("N42").Value = 22000
$ws.Range("O42").Value = 24000
$ws.Range("P42").Value = 23000
$ws.Range("Q42").Value = "`$/caja 16 kilos"
$ws.Range("R42").Value = "Paine"
$ws.Range("S42").Value = 1438
$ws.Range("T42").Value = 16

$ws.Range("A43").Value = 6
$ws.Range("B43").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C43").Value = "Metropolitana"
$ws.Range("D43").Value = 44168
$ws.Range("E43").Value = 13
$ws.Range("F43").Value = "Fruta"
$ws.Range("G43").Value = 100103
$ws.Range("H43").Value = "Frutos de hueso (carozo)"
$ws.Range("I43").Value = 100103003
$ws.Range("J43").Value = "Damasco"
$ws.Range("K43").Value = "Castle Brite"
$ws.Range("L43").Value = "Primera"
$ws.Range("M43").Value = 120
$ws.Range("N43").Value = 18000
$ws.Range("O43").Value = 18000
$ws.Range("P43").Value = 18000
$ws.Range("Q43").Value = "`$/caja 18 kilos"
$ws.Range("R43").Value = "Región Metropolitana"
$ws.Range("S43").Value = 1000
$ws.Range("T43").Value = 18

$ws.Range("A44").Value = 6
$ws.Range("B44").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C44").Value = "Metropolitana"
$ws.Range("D44").Value = 44168
$ws.Range("E44").Value = 13
$ws.Range("F44").Value = "Fruta"
$ws.Range("G44").Value = 100103
$ws.Range("H44").Value = "Frutos de hueso (carozo)"
$ws.Range("I44").Value = 100103003
$ws.Range("J44").Value = "Damasco"
$ws.Range("K44").Value = "Dina"
$ws.Range("L44").Value = "Primera"
$ws.Range("M44").Value = 10
$ws.Range("N44").Value = 500000
$ws.Range("O44").Value = 500000
$ws.Range("P44").Value = 500000
$ws.Range("Q44").Value = "`$/bins (500 kilos)"
$ws.Range("R44").Value = "Región Metropolitana"
$ws.Range("S44").Value = 1000
$ws.Range("T44").Value = 500

$ws.Range("A45").Value = 6
$ws.Range("B45").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C45").Value = "Metropolitana"
$ws.Range("D45").Value = 44166
$ws.Range("E45").Value = 13
$ws.Range("F45").Value = "Fruta"
$ws.Range("G45").Value = 100103
$ws.Range("H45").Value = "Frutos de hueso (carozo)"
$ws.Range("I45").Value = 100103003
$ws.Range("J45").Value = "Damasco"
$ws.Range("K45").Value = "Dina"
$ws.Range("L45").Value = "Primera"
$ws.Range("M45").Value = 150
$ws.Range("N45").Value = 13000
$ws.Range("O45").Value = 13000
$ws.Range("P45").Value = 13000
$ws.Range("Q45").Value = "`$/bandeja 6 kilos"
$ws.Range("R45").Value = "Región Metropolitana"
$ws.Range("S45").Value = 2167
$ws.Range("T45").Value = 6

$ws.Range("A46").Value = 6
$ws.Range("B46").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C46").Value = "Metropolitana"
$ws.Range("D46").Value = 44166
$ws.Range("E46").Value = 13
$ws.Range("F46").Value = "Fruta"
$ws.Range("G46").Value = 100103
$ws.Range("H46").Value = "Frutos de hueso (carozo)"
$ws.Range("I46").Value = 100103003
$ws.Range("J46").Value = "Damasco"
$ws.Range("K46").Value = "Dina"
$ws.Range("L46").Value = "Segunda"
$ws.Range("M46").Value = 150
$ws.Range("N46").Value = 11000
$ws.Range("O46").Value = 11000
$ws.Range("P46").Value = 11000
$ws.Range("Q46").Value = "`$/bandeja 6 kilos"
$ws.Range("R46").Value = "Región Metropolitana"
$ws.Range("S46").Value = 1833
$ws.Range("T46").Value = 6

$ws.Range("A47").Value = 6
$ws.Range("B47").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C47").Value = "Metropolitana"
$ws.Range("D47").Value = 44525
$ws.Range("E47").Value = 13
$ws.Range("F47").Value = "Fruta"
$ws.Range("G47").Value = 100103
$ws.Range("H47").Value = "Frutos de hueso (carozo)"
$ws.Range("I47").Value = 100103003
$ws.Range("J47").Value = "Damasco"
$ws.Range("K47").Value = "Castle Brite"
$ws.Range("L47").Value = "Primera"
$ws.Range("M47").Value = 12
$ws.Range("N47").Value = 580000
$ws.Range("O47").Value = 580000
$ws.Range("P47").Value = 580000
$ws.Range("Q47").Value = "`$/bins (500 kilos)"
$ws.Range("R47").Value = "Paine"
$ws.Range("S47").Value = 1160
$ws.Range("T47").Value = 500

$ws.Range("A48").Value = 6
$ws.Range("B48").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C48").Value = "Metropolitana"
$ws.Range("D48").Value = 44160
$ws.Range("E48").Value = 13
$ws.Range("F48").Value = "Fruta"
$ws.Range("G48").Value = 100103
$ws.Range("H48").Value = "Frutos de hueso (carozo)"
$ws.Range("I48").Value = 100103003
$ws.Range("J48").Value = "Damasco"
$ws.Range("K48").Value = "Dina"
$ws.Range("L48").Value = "Especial"
$ws.Range("M48").Value = 50
$ws.Range("N48").Value = 16000
$ws.Range("O48").Value = 16000
$ws.Range("P48").Value = 16000
$ws.Range("Q48").Value = "`$/caja 16 kilos"
$ws.Range("R48").Value = "Región Metropolitana"
$ws.Range("S48").Value = 1000
$ws.Range("T48").Value = 16

$ws.Range("A49").Value = 6
$ws.Range("B49").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C49").Value = "Metropolitana"
$ws.Range("D49").Value = 44160
$ws.Range("E49").Value = 13
$ws.Range("F49").Value = "Fruta"
$ws.Range("G49").Value = 100103
$ws.Range("H49").Value = "Frutos de hueso (carozo)"
$ws.Range("I49").Value = 100103003
$ws.Range("J49").Value = "Damasco"
$ws.Range("K49").Value = "Dina"
$ws.Range("L49").Value = "Primera"
$ws.Range("M49").Value = 270
$ws.Range("N49").Value = 14000
$ws.Range("O49").Value = 14000
$ws.Range("P49").Value = 14000
$ws.Range("Q49").Value = "`$/caja 15 kilos"
$ws.Range("R49").Value = "Región Metropolitana"
$ws.Range("S49").Value = 933
$ws.Range("T49").Value = 15

$ws.Range("A50").Value = 6
$ws.Range("B50").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C50").Value = "Metropolitana"
$ws.Range("D50").Value = 44160
$ws.Range("E50").Value = 13
$ws.Range("F50").Value = "Fruta"
$ws.Range("G50").Value = 100103
$ws.Range("H50").Value = "Frutos de hueso (carozo)"
$ws.Range("I50").Value = 100103003
$ws.Range("J50").Value = "Damasco"
$ws.Range("K50").Value = "Dina"
$ws.Range("L50").Value = "Primera"
$ws.Range("M50").Value = 70
$ws.Range("N50").Value = 14000
$ws.Range("O50").Value = 14000
$ws.Range("P50").Value = 14000
$ws.Range("Q50").Value = "`$/caja 16 kilos"
$ws.Range("R50").Value = "Región Metropolitana"
$ws.Range("S50").Value = 875
$ws.Range("T50").Value = 16

$ws.Range("A51").Value = 6
$ws.Range("B51").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C51").Value = "Metropolitana"
$ws.Range("D51").Value = 44160
$ws.Range("E51").Value = 13
$ws.Range("F51").Value = "Fruta"
$ws.Range("G51").Value = 100103
$ws.Range("H51").Value = "Frutos de hueso (carozo)"
$ws.Range("I51").Value = 100103003
$ws.Range("J51").Value = "Damasco"
$ws.Range("K51").Value = "Dina"
$ws.Range("L51").Value = "Segunda"
$ws.Range("M51").Value = 30
$ws.Range("N51").Value = 12000
$ws.Range("O51").Value = 12000
$ws.Range("P51").Value = 12000
$ws.Range("Q51").Value = "`$/caja 16 kilos"
$ws.Range("R51").Value = "Región Metropolitana"
$ws.Range("S51").Value = 750
$ws.Range("T51").Value = 16

$ws.Range("A52").Value = 6
$ws.Range("B52").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C52").Value = "Metropolitana"
$ws.Range("D52").Value = 44160
$ws.Range("E52").Value = 13
$ws.Range("F52").Value = "Fruta"
$ws.Range("G52").Value = 100103
$ws.Range("H52").Value = "Frutos de hueso (carozo)"
$ws.Range("I52").Value = 100103003
$ws.Range("J52").Value = "Damasco"
$ws.Range("K52").Value = "Dina"
$ws.Range("L52").Value = "Tercera"
$ws.Range("M52").Value = 30
$ws.Range("N52").Value = 10000
$ws.Range("O52").Value = 10000
$ws.Range("P52").Value = 10000
$ws.Range("Q52").Value = "`$/caja 16 kilos"
$ws.Range("R52").Value = "Región Metropolitana"
$ws.Range("S52").Value = 625
$ws.Range("T52").Value = 16

$ws.Range("A53").Value = 6
$ws.Range("B53").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C53").Value = "Metropolitana"
$ws.Range("D53").Value = 44181
$ws.Range("E53").Value = 13
$ws.Range("F53").Value = "Fruta"
$ws.Range("G53").Value = 100103
$ws.Range("H53").Value = "Frutos de hueso (carozo)"
$ws.Range("I53").Value = 100103003
$ws.Range("J53").Value = "Damasco"
$ws.Range("K53").Value = "Dina"
$ws.Range("L53").Value = "Primera"
$ws.Range("M53").Value = 115
$ws.Range("N53").Value = 11000
$ws.Range("O53").Value = 11000
$ws.Range("P53").Value = 11000
$ws.Range("Q53").Value = "`$/caja 16 kilos"
$ws.Range("R53").Value = "Región Metropolitana"
$ws.Range("S53").Value = 688
$ws.Range("T53").Value = 16

$ws.Range("A54").Value = 6
$ws.Range("B54").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C54").Value = "Metropolitana"
$ws.Range("D54").Value = 44181
$ws.Range("E54").Value = 13
$ws.Range("F54").Value = "Fruta"
$ws.Range("G54").Value = 100103
$ws.Range("H54").Value = "Frutos de hueso (carozo)"
$ws.Range("I54").Value = 100103003
$ws.Range("J54").Value = "Damasco"
$ws.Range("K54").Value = "Dina"
$ws.Range("L54").Value = "Segunda"
$ws.Range("M54").Value = 70
$ws.Range("N54").Value = 9000
$ws.Range("O54").Value = 9000
$ws.Range("P54").Value = 9000
$ws.Range("Q54").Value = "`$/caja 16 kilos"
$ws.Range("R54").Value = "Región Metropolitana"
$ws.Range("S54").Value = 562
$ws.Range("T54").Value = 16

$ws.Range("A55").Value = 6
$ws.Range("B55").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C55").Value = "Metropolitana"
$ws.Range("D55").Value = 44181
$ws.Range("E55").Value = 13
$ws.Range("F55").Value = "Fruta"
$ws.Range("G55").Value = 100103
$ws.Range("H55").Value = "Frutos de hueso (carozo)"
$ws.Range("I55").Value = 100103003
$ws.Range("J55").Value = "Damasco"
$ws.Range("K55").Value = "Modesto"
$ws.Range("L55").Value = "Primera"
$ws.Range("M55").Value = 12
$ws.Range("N55").Value = 450000
$ws.Range("O55").Value = 450000
$ws.Range("P55").Value = 450000
$ws.Range("Q55").Value = "`$/bins (500 kilos)"
$ws.Range("R55").Value = "Región Metropolitana"
$ws.Range("S55").Value = 900
$ws.Range("T55").Value = 500

$ws.Range("A56").Value = 6
$ws.Range("B56").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C56").Value = "Metropolitana"
$ws.Range("D56").Value = 44181
$ws.Range("E56").Value = 13
$ws.Range("F56").Value = "Fruta"
$ws.Range("G56").Value = 100103
$ws.Range("H56").Value = "Frutos de hueso (carozo)"
$ws.Range("I56").Value = 100103003
$ws.Range("J56").Value = "Damasco"
$ws.Range("K56").Value = "Modesto"
$ws.Range("L56").Value = "Segunda"
$ws.Range("M56").Value = 6
$ws.Range("N56").Value = 380000
$ws.Range("O56").Value = 380000
$ws.Range("P56").Value = 380000
$ws.Range("Q56").Value = "`$/bins (500 kilos)"
$ws.Range("R56").Value = "Región Metropolitana"
$ws.Range("S56").Value = 760
$ws.Range("T56").Value = 500

$ws.Range("A57").Value = 6
$ws.Range("B57").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C57").Value = "Metropolitana"
$ws.Range("D57").Value = 44189
$ws.Range("E57").Value = 13
$ws.Range("F57").Value = "Fruta"
$ws.Range("G57").Value = 100103
$ws.Range("H57").Value = "Frutos de hueso (carozo)"
$ws.Range("I57").Value = 100103003
$ws.Range("J57").Value = "Damasco"
$ws.Range("K57").Value = "Dina"
$ws.Range("L57").Value = "Especial"
$ws.Range("M57").Value = 45
$ws.Range("N57").Value = 18000
$ws.Range("O57").Value = 18000
$ws.Range("P57").Value = 18000
$ws.Range("Q57").Value = "`$/caja 18 kilos"
$ws.Range("R57").Value = "Región de O'Higgins"
$ws.Range("S57").Value = 1000
$ws.Range("T57").Value = 18
$ws.Range("D57").NumberFormat = "YYYY-MM-DD HH:MM:SS"

$ws.Range("A58").Value = 6
$ws.Range("B58").Value = "Mercado Mayorista Lo Valledor de Santiago"
$ws.Range("C58").Value = "Metropolitana"
$ws.Range("D58").Value = 44189
$ws.Range("E58").Value = 13
$ws.Range("F58").Value = "Fruta"
$ws.Range("G58").Value = 100103
$ws.Range("H58").Value = "Frutos de hueso (carozo)"
$ws.Range("I58").Value = 100103003
$ws.Range("J58").Value = "Damasco"
$ws.Range("K58").Value = "Dina"
$ws.Range("L58").Value = "Primera"
$ws.Range("M58").Value = 50
$ws.Range("N58").Value = 16000
$ws.Range("O58").Value = 16000
$ws.Range("P58").Value = 16000
$ws.Range("Q58").Value = "`$/caja 18 kilos"
$ws.Range("R58").Value = "Región de O'Higgins"
$ws.Range("S58").Value = 889
$ws.Range("T58").Value = 18
$ws.Range("D58").NumberFormat = "YYYY-MM-DD HH:MM:SS"
